# Rename the "Effet à obtenir" sheet to "#Effet à obtenir".
# Renaming automatically updates formulas that reference the sheet
# (e.g. C21/D21 on the #Sommaire sheet) while leaving hyperlink
# "location" targets untouched, matching native Excel behaviour.
$wb = $excel.ActiveWorkbook
$target = $wb.Worksheets.Item("Effet à obtenir")
$target.Name = "#Effet à obtenir"

# Activate the renamed sheet and move the selection, which is what
# flips tabSelected from the previously active sheet onto this one
# and records the new selection/active-cell state.
$target.Activate()
$target.Range("K27").Select()
